# Rename the two worksheets:
#   "Sheet 1" -> "Personal_Info"
#   "Sheet 2" -> "Product_Info"
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    if ($ws.Name -eq "Sheet 1") {
        $ws.Name = "Personal_Info"
    } elseif ($ws.Name -eq "Sheet 2") {
        $ws.Name = "Product_Info"
    }
}
